$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO CMS")

# Mark "Löschen von Bestellungen" (row 2) as done: copy the "done" formatting
# (green fill / white bold text, style used by rows 3 and 9) onto B2, then
# update its text from "offen" to "done".
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B2").Value2 = "done"

# Add a new TODO row for the product editor: moving pre-orders when the
# amount changes.
$ws.Range("A12").Value2 = "Vorbestellungen verschieben, wenn Amount geändert wird"
$ws.Range("B4").Copy()
$ws.Range("B12").PasteSpecial(-4122)
$ws.Range("B12").Value2 = "offen"

# Update the active selection shown when the workbook was last saved.
$ws.Activate()
$ws.Range("C5").Select()
